$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the existing header formatting (bold, centered, bordered) from A1
# instead of inventing a new style entry.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every team row gets the same season record (91-71-0).
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 91
    $ws.Cells.Item($row, 31).Value = 71
    $ws.Cells.Item($row, 32).Value = 0
}
